$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
# Insert all 4 new rows first (top to bottom so positions don't shift)
$ws.Rows("29:29").Insert(-4121)
$ws.Rows("45:45").Insert(-4121)
$ws.Rows("51:51").Insert(-4121)
$ws.Rows("67:67").Insert(-4121)

# Row 29: A29 needs style like row27 col A/B (s=57) but without alignment (target s=58)
# We'll copy format from A68 B68 (s46/s47 group) then adjust - but actually target uses fontId=9 fillId=3 xfId=58 (no alignment)
# The closest existing reference cell with fontId9 fillId3 is A27/B27 (s=57) but that has alignment vertical=center.
# Since no exact match exists pre-built, we copy from A27 (closest) then it will likely reuse s=57 (with alignment) - not exact match of target s=58.
# Instead, build via direct formatting: font9 (bold, black), fillId=3 (light blue/grey). Let's inspect fonts/fills to replicate with Interior/Font properties, and xfId base style "58" (different base reference), which Excel cannot set directly via COM normally - it's tied to "cell style" application.
Write-Host "rows inserted"
